$d = $word.ActiveDocument

# Locate the paragraph that begins "Trang web đồng thời cũng ..." (the phrase
# spans the boundary between the "Trang web " run and the "đồng thời cũng " run,
# so searching for it pins down the exact split point unambiguously).
$find = $d.Content
$find.Find.Text = "Trang web đồng thời cũng"
$find.Find.Forward = $true
$find.Find.MatchCase = $true
$find.Find.Execute() | Out-Null
$paraStart = $find.Start
$splitPos = $paraStart + 10   # length of "Trang web " (including trailing space)

# Figure out which paragraph index that is, so we can re-address it reliably
# after the split shifts character offsets around.
$count = $d.Paragraphs.Count
$paraIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $paraStart) {
        $paraIndex = $i
        break
    }
}

# Split the paragraph right before "đồng thời cũng" so "Trang web " is left
# alone as the sole content of its own paragraph. The original paragraph
# mark/properties stay with this first half; a brand-new (bare) paragraph is
# created to hold the remainder ("đồng thời cũng ứng dụng ...").
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertParagraphBefore()

# After the split: paragraph $paraIndex == "Trang web " (+ its paragraph mark);
# paragraph $paraIndex + 1 == the new bare paragraph starting with
# "đồng thời cũng ...".
$firstPara = $d.Paragraphs.Item($paraIndex)
$secondPara = $d.Paragraphs.Item($paraIndex + 1)

# Grab a formatted copy of "Trang web " (keeps its run formatting/properties)
# before we overwrite it, excluding the trailing paragraph mark.
$trangWebRng = $d.Range($firstPara.Range.Start, $firstPara.Range.End - 1)
$formattedTrangWeb = $trangWebRng.FormattedText

# Re-insert the formatted "Trang web " text, unmodified, at the very top of the
# new second paragraph (right before "đồng thời cũng").
$insRng = $d.Range($secondPara.Range.Start, $secondPara.Range.Start)
$insRng.FormattedText = $formattedTrangWeb

# Finally, turn the original "Trang web " run (now the sole content of the
# first half of the split) into the new sentence about the caro game.
$trangWebRng2 = $d.Range($firstPara.Range.Start, $firstPara.Range.End - 1)
$trangWebRng2.Text = "Ngoài ra, bạn có thể kết nối với một người dùng khác để chơi caro ngay trên trang web."
